# Insert a new data row right before the current row 258 (shifts 258..369
# down to 259..370, extending the used range to A1:R370), then populate the
# newly-inserted row 258 with its own record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(258).Insert()

$ws.Range("A258").Value = 9
$ws.Range("B258").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C258").Value = "Metropolitana"
$ws.Range("D258").Value = 44636
$ws.Range("E258").Value = 13
$ws.Range("F258").Value = 100112039
$ws.Range("G258").Value = "Ciboulette"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 450
$ws.Range("K258").Value = 1500
$ws.Range("L258").Value = 1800
$ws.Range("M258").Value = 1633
$ws.Range("N258").Value = "$/docena de atados"
$ws.Range("O258").Value = "Provincia de Chacabuco"
$ws.Range("P258").Value = 544
$ws.Range("Q258").Value = 3
$ws.Range("R258").Value = "Hortaliza"
